# "Generate Report for Handback" - refresh the handoff/handback timestamps
# recorded for the a58f0387-6e9c-4f5f-8ce6-b8cafc7a37d6 item (which shares
# its timestamp strings with the a93bb8cf-a4b8-423a-827f-05bac5951e1e item)
# on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-11 22:14:52"
$zhcn.Range("E4").Value = "2016-03-11 22:14:52"
$zhcn.Range("H3").Value = "2016-03-11 22:15:12"
$zhcn.Range("H4").Value = "2016-03-11 22:15:12"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-11 22:14:56"
$dede.Range("E4").Value = "2016-03-11 22:14:56"
$dede.Range("H3").Value = "2016-03-11 22:15:12"
$dede.Range("H4").Value = "2016-03-11 22:15:12"
